# Fix the "PRC" -> "RPC" typo on slide 14 ("4. Service Communication Design"),
# in the "Content Placeholder 2" shape's first bullet:
#   "... HTTP, AMQP, and PRC. "  ->  "... HTTP, AMQP, and RPC. "
#
# Doing this through TextRange.Characters(start, length) (rather than
# rewriting the whole TextRange.Text) mirrors how PowerPoint itself records
# an in-place retype: it only touches the "and PRC" substring, which is
# exactly what keeps the run split boundaries ("...AMQP, " | "and RPC" | ". ")
# matching the authored edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(14)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$full = $tr.Text
$start = $full.IndexOf("and PRC") + 1   # COM Characters() is 1-indexed

$chars = $tr.Characters($start, 7)
$chars.Text = "and RPC"
